$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header values (number of repetitions / identifiers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON values
$ws.Range("B2").Value = 62.077818775276398
$ws.Range("C2").Value = 49.138467275708635
$ws.Range("D2").Value = 65.128684398943349
$ws.Range("E2").Value = 51.940719555127188

# Row 3 - STR values
$ws.Range("B3").Value = 63.589254376254942
$ws.Range("C3").Value = 44.929400036024902
$ws.Range("D3").Value = 72.762927595038093
$ws.Range("E3").Value = 48.724044589012166

# Update the active selection to match the updated data range
[void]$ws.Range("B1:E3").Select()
